# Test Herjedal.xlsx edit
#
# The sheet used to contain an old/duplicate EBITDA block (absolute values in
# row 11 + a % growth/margin row in row 12) sitting between the "% vekst"
# row (10) and the "EBITDA" / "EBITDA (pre IFRS)" / "% margin" block that used
# to live in rows 13-15. That stale block is removed, which shifts the
# EBITDA block up so it now sits directly under the revenue growth row.
#
# Also: the leftover "-" placeholder in B10 is cleared, the growth/margin
# percentage cells get a new custom number format (0.0 %), and the saved
# selection/active cell moves to F26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old absolute-EBITDA row (11) and its % growth/margin row (12).
# Rows 13:15 (EBITDA / EBITDA (pre IFRS) / % margin) shift up to 11:13 and
# their formulas are automatically re-pointed (B12/B9 etc, instead of B14/B9).
$ws.Rows("11:12").Delete()

# The "-" text that used to sit next to "% vekst" is no longer needed.
$ws.Range("B10").ClearContents()

# Give the growth-rate / margin formula cells the new "0.0 %" number format.
$ws.Range("C10:G10").NumberFormat = "0.0\ %"
$ws.Range("B13:G13").NumberFormat = "0.0\ %"

# Match the saved cursor position/selection in the workbook view.
$ws.Range("F26").Select()
